# Turn keyword from absetn to abs.
# Also inserts a new "Groupe" column (value "1A" for every player) on the
# first three sheets (Fitness, Basket M, Basket F).

$wb = $excel.ActiveWorkbook

$sheetNames = @("Fitness", "Basket M", "Basket F")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Insert a new first column ("Groupe") and shift everything else right.
    $ws.Columns.Item(1).EntireColumn.Insert()

    $ws.Range("A1").Value = "Groupe"
    $ws.Range("A2:A6").Value = "1A"
}

# On the "Basket F" sheet, the absence keyword changes from "Absent" to "Abs.".
$wsF = $wb.Worksheets.Item("Basket F")
$used = $wsF.Range("E2:R6")
for ($r = 1; $r -le $used.Rows.Count; $r++) {
    for ($c = 1; $c -le $used.Columns.Count; $c++) {
        $cell = $used.Cells.Item($r, $c)
        if ($cell.Value -eq "Absent") {
            $cell.Value = "Abs."
        }
    }
}

# Selections left by the editing session.
$wsFitness = $wb.Worksheets.Item("Fitness")
$wsFitness.Range("A1:A1048576").Select()
$wsFitness.Application.ActiveWindow.Zoom = 86

$wsBasketM = $wb.Worksheets.Item("Basket M")
$wsBasketM.Range("A1:A1048576").Select()

$wsF.Range("I10").Select()
$wsF.Activate()
